$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.357.36"
$ws.Range("E2").Value = "  +1.10%  "

$ws.Range("D3").Value = "1.780.52"
$ws.Range("E3").Value = "  +3.71%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "'314.07"
$ws.Range("E5").Value = "  +1.07%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").Value = "'0.5252"
$ws.Range("E7").Value = "  +9.73%  "

$ws.Range("D8").Value = "'0.3741"
$ws.Range("E8").Value = "  +8.22%  "

$ws.Range("D9").Value = "'42.77"
$ws.Range("E9").Value = "  +1.51%  "

$ws.Range("D10").Value = "'0.07406"
$ws.Range("E10").Value = "  +2.29%  "

$ws.Range("D11").Value = "'1.095"
$ws.Range("E11").Value = "  +5.28%  "

$ws.Range("D12").Value = "'0.9997"
$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("E13").Value = "  +5.31%  "

$ws.Range("D14").Value = "'6.120"
$ws.Range("E14").Value = "  +4.74%  "

$ws.Range("D15").Value = "1.778.44"
$ws.Range("E15").Value = "  +3.63%  "

$ws.Range("D16").Value = "'6.987"
$ws.Range("E16").Value = "  +2.52%  "

$ws.Range("D17").Value = "'89.67"
$ws.Range("E17").Value = "  +2.64%  "

$ws.Range("D18").Value = "'0.00001058"
$ws.Range("E18").Value = "  +2.46%  "

$ws.Range("D19").Value = "'0.06440"
$ws.Range("E19").Value = "  +0.73%  "

$ws.Range("D20").Value = "'0.9999"
$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").Value = "'16.79"
$ws.Range("E21").Value = "  +2.16%  "

$ws.Range("D22").Value = "'5.895"
$ws.Range("E22").Value = "  +5.10%  "

$ws.Range("D23").Value = "27.391.38"
$ws.Range("E23").Value = "  +0.97%  "

$ws.Range("E24").Value = "  +4.24%  "

$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("D26").Value = "'155.60"
$ws.Range("E26").Value = "  +3.64%  "

$ws.Range("E27").Value = "  +1.50%  "

$ws.Range("D28").Value = "'2.362"
$ws.Range("E28").Value = "  +13.46%  "

$ws.Range("D29").Value = "1.985.27"
$ws.Range("E29").Value = "  +3.61%  "

$ws.Range("D30").Value = "'121.21"
$ws.Range("E30").Value = "  +0.33%  "

$ws.Range("D31").Value = "'1.094"
$ws.Range("E31").Value = "  +5.54%  "

$ws.Range("D32").Value = "'0.1013"
$ws.Range("E32").Value = "  +10.04%  "

$ws.Range("D33").Value = "'5.599"
$ws.Range("E33").Value = "  +5.39%  "

$ws.Range("D34").Value = "'3.625"
$ws.Range("E34").Value = "  +0.64%  "

$ws.Range("D35").Value = "'0.02261"
$ws.Range("E35").Value = "  +3.56%  "

$ws.Range("D36").Value = "'0.05985"
$ws.Range("E36").Value = "  +2.46%  "

$ws.Range("D37").Value = "'11.34"
$ws.Range("E37").Value = "  +3.80%  "

$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2051"
$ws.Range("E38").Value = "  +3.40%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'4.892"
$ws.Range("E39").Value = "  +3.74%  "

$ws.Range("D40").Value = "'0.6129"
$ws.Range("E40").Value = "  +3.23%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'8.217"
$ws.Range("E41").Value = "  +9.70%  "

$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'1.433"
$ws.Range("E42").Value = "  -2.68%  "

$ws.Range("E43").Value = "  +4.41%  "

$ws.Range("E44").Value = "  +3.45%  "

$ws.Range("D45").Value = "'0.5793"
$ws.Range("E45").Value = "  +4.27%  "

$ws.Range("D46").Value = "'3.623"
$ws.Range("E46").Value = "  +0.97%  "

$ws.Range("D47").Value = "'121.69"
$ws.Range("E47").Value = "  +2.47%  "

$ws.Range("D48").Value = "'1.897"
$ws.Range("E48").Value = "  +4.21%  "

$ws.Range("D49").Value = "'1.115"
$ws.Range("E49").Value = "  +2.49%  "

$ws.Range("E50").Value = "  +0.78%  "

$ws.Range("D51").Value = "'70.95"
$ws.Range("E51").Value = "  +2.45%  "
